# CRUD-Division Address.xlsx update
# - Populate CRUD testcases (Feature category) on both sheets (rows 3-5 on
#   Create_DivisionAddress, rows 2-5 on Edit_DivisionAddress incl. VAT Class
#   (City) variations: New Mumbai / Pune / Delhi).
# - Re-create the Website / Contact Email hyperlinks for the new rows.
# - Re-point the active sheet / selections to match the new working state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Create_DivisionAddress")
$ws2 = $wb.Worksheets.Item("Edit_DivisionAddress")

# ---------------------------------------------------------------------
# Sheet 1: Create_DivisionAddress - fill in rows 3-5 with the same
# address data pattern already present in row 2.
# ---------------------------------------------------------------------
$rows1 = 3,4,5
foreach ($r in $rows1) {
    $ws1.Range("B$r").Value = "Div Street"
    $ws1.Range("C$r").Value = "Mumbai"
    $ws1.Range("D$r").Value = 400042
    $ws1.Range("E$r").Value = 111000
    $ws1.Range("F$r").Value = "www.google.com"
    $ws1.Range("H$r").Value = "MH"
    $ws1.Range("I$r").Value = "India"
    $ws1.Range("J$r").Value = 121212
    $ws1.Range("K$r").Value = "pkakade@rootstock.com"
}

foreach ($r in $rows1) {
    $ws1.Hyperlinks.Add($ws1.Range("F$r"), "http://www.google.com/") | Out-Null
}
foreach ($r in $rows1) {
    $ws1.Hyperlinks.Add($ws1.Range("K$r"), "mailto:pkakade@rootstock.com") | Out-Null
}
foreach ($r in $rows1) {
    $ws1.Range("F$r").Style = "Hyperlink"
    $ws1.Range("K$r").Style = "Hyperlink"
}

$ws1.Columns.Item(1).AutoFit() | Out-Null
$ws1.Columns.Item(1).ColumnWidth = 14.59
$ws1.Columns.Item(8).AutoFit() | Out-Null
$ws1.Columns.Item(8).ColumnWidth = 12.26

# ---------------------------------------------------------------------
# Sheet 2: Edit_DivisionAddress - build out the full header row plus
# four data rows (City varies: New Mumbai / Pune / Delhi / Pune).
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "AddressType"
$ws2.Range("B1").Value = "Street"
$ws2.Range("C1").Value = "City"
$ws2.Range("D1").Value = "Zip/PostalCode"
$ws2.Range("E1").Value = "Phone"
$ws2.Range("F1").Value = "Website"
$ws2.Range("G1").Value = "Shipping Zone"
$ws2.Range("H1").Value = "State/Province"
$ws2.Range("I1").Value = "Country"
$ws2.Range("J1").Value = "Fax"
$ws2.Range("K1").Value = "Contact Email"

$ws2.Range("A2").Value = "All"
$ws2.Range("A3").Value = "Receiving Address"
$ws2.Range("A4").Value = "Shipping Address"
$ws2.Range("A5").Value = "Remit Address"

$cities = @{ 2 = "New Mumbai"; 3 = "Pune"; 4 = "Delhi"; 5 = "Pune" }
foreach ($r in 2,3,4,5) {
    $ws2.Range("B$r").Value = "Div Street"
    $ws2.Range("C$r").Value = $cities[$r]
    $ws2.Range("D$r").Value = 400042
    $ws2.Range("E$r").Value = 111000
    $ws2.Range("F$r").Value = "www.google.com"
    $ws2.Range("H$r").Value = "MH"
    $ws2.Range("I$r").Value = "India"
    $ws2.Range("J$r").Value = 121212
    $ws2.Range("K$r").Value = "pkakade@rootstock.com"
}

$ws2.Hyperlinks.Add($ws2.Range("F2"), "http://www.google.com/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("K2"), "mailto:pkakade@rootstock.com") | Out-Null
foreach ($r in 3,4,5) {
    $ws2.Hyperlinks.Add($ws2.Range("F$r"), "http://www.google.com/") | Out-Null
}
foreach ($r in 3,4,5) {
    $ws2.Hyperlinks.Add($ws2.Range("K$r"), "mailto:pkakade@rootstock.com") | Out-Null
}
foreach ($r in 2,3,4,5) {
    $ws2.Range("F$r").Style = "Hyperlink"
    $ws2.Range("K$r").Style = "Hyperlink"
}

$ws2.Columns.Item(11).AutoFit() | Out-Null
$ws2.Columns.Item(11).ColumnWidth = 20.92

# ---------------------------------------------------------------------
# Selections / active sheet - match the new working state: the editor
# left Create_DivisionAddress with D29 selected, and switched to
# Edit_DivisionAddress (whole data range selected) as the active tab.
# ---------------------------------------------------------------------
$ws1.Range("D29").Select() | Out-Null
$ws2.Range("A1:K5").Select() | Out-Null
$ws2.Activate()
